$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 96 (pushes the former rows 96-98 down to 97-99)
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row with this week's price entry
$ws.Range("A96").Value2 = 10
$ws.Range("B96").Value2 = "Vega Modelo de Temuco"
$ws.Range("C96").Value2 = "La Araucanía"
$ws.Range("D96").Value2 = 44753
$ws.Range("E96").Value2 = 9
$ws.Range("F96").Value2 = 100114002
$ws.Range("G96").Value2 = "Camote"
$ws.Range("H96").Value2 = "Sin especificar"
$ws.Range("I96").Value2 = "Primera"
$ws.Range("J96").Value2 = 80
$ws.Range("K96").Value2 = 20000
$ws.Range("L96").Value2 = 20000
$ws.Range("M96").Value2 = 20000
$ws.Range("N96").Value2 = "`$/malla 20 kilos"
$ws.Range("O96").Value2 = "Perú"
$ws.Range("P96").Value2 = 1000
$ws.Range("Q96").Value2 = 20
$ws.Range("R96").Value2 = "Hortaliza"
